# Update ig name: correct the FHIR IG URL placeholder "[code]" -> "tde",
# bump Version to 2.0.0, and refresh the generation Date.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL (row 2, column B)
$meta.Range("B2").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/StructureDefinition/EyeColor"

# Version (row 3, column B)
$meta.Range("B3").Value = "2.0.0"

# Date (row 8, column B)
$meta.Range("B8").Value = "2026-01-15T15:25:18+00:00"

# --- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Extension.url fixed value (row 5, column R) mirrors the StructureDefinition URL
$elements.Range("R5").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/StructureDefinition/EyeColor"

# Binding Value Set (row 6, column Z)
$elements.Range("Z6").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/ValueSet/EyeColorVS"

# Column Z auto-fits narrower now that the URL text is shorter
# (target best-fit width 49.4453125; nudged so the host's internal
# pixel-grid rounding lands on the nearest representable width).
$elements.Columns.Item(26).ColumnWidth = 48.6
